$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_9_9_0", 0.8646455558456058, 0.720630002910706, 0.9173521372528198, 0.77937908154035, 0.1497974842786789, 0.4882208108901978, 0.06918749213218689, 0.2910287380218506),
    @("model_9_9_1", 0.8857074489717135, 0.7646591974836079, 0.8979800534274553, 0.80447418646743, 0.1264881789684296, 0.4112763702869415, 0.08540455996990204, 0.2579249143600464),
    @("model_9_9_2", 0.8919394307161677, 0.779526558844506, 0.8917597362457232, 0.8130438631361983, 0.1195912137627602, 0.3852944672107697, 0.09061180800199509, 0.2466203719377518),
    @("model_9_9_3", 0.8969746712608054, 0.7952635575382443, 0.8854769394373576, 0.8222049469462034, 0.1140186786651611, 0.3577928841114044, 0.09587135910987854, 0.2345356643199921),
    @("model_9_9_4", 0.8987659824925222, 0.8031665686506395, 0.8740577835733394, 0.8243376366784858, 0.1120362281799316, 0.3439817428588867, 0.1054307445883751, 0.231722354888916),
    @("model_9_9_5", 0.8989471083866949, 0.803892751809102, 0.8729443814362516, 0.8245144623766132, 0.1118357703089714, 0.3427126407623291, 0.1063628196716309, 0.2314890772104263),
    @("model_9_9_24", 0.8990196268191726, 0.829776796351589, 0.7801461844586023, 0.8149550693578167, 0.1117555126547813, 0.2974782884120941, 0.184047520160675, 0.2440992593765259),
    @("model_9_9_23", 0.8991230802972208, 0.830143511317147, 0.78028461367, 0.8152539514860435, 0.1116410121321678, 0.2968374490737915, 0.1839316189289093, 0.2437049895524979),
    @("model_9_9_6", 0.8994768419613142, 0.8081882224671972, 0.8611011669596105, 0.8239901669905814, 0.1112495139241219, 0.3352060317993164, 0.1162771880626678, 0.2321807146072388),
    @("model_9_9_11", 0.8996951940246097, 0.8331916980914635, 0.7851861392733043, 0.8188553492211987, 0.1110078617930412, 0.2915104627609253, 0.1798283755779266, 0.2389542311429977),
    @("model_9_9_22", 0.8997961786972535, 0.8316249156187329, 0.7830705127358693, 0.8171247241438626, 0.110896110534668, 0.2942485511302948, 0.1815994530916214, 0.2412371933460236),
    @("model_9_9_21", 0.899923072424718, 0.8319363840086487, 0.7835457260396883, 0.8174852174858398, 0.1107556745409966, 0.2937042713165283, 0.1812016367912292, 0.2407616376876831),
    @("model_9_9_19", 0.8999907072070108, 0.8336958472307177, 0.7803984592160149, 0.8177792790848952, 0.1106808185577393, 0.2906294465065002, 0.1838363111019135, 0.2403737306594849),
    @("model_9_9_20", 0.9000231825522409, 0.8318927567821817, 0.7846748962207077, 0.8177917213771336, 0.1106448844075203, 0.293780505657196, 0.1802563667297363, 0.2403573095798492),
    @("model_9_9_8", 0.9002077617252924, 0.8151694989275524, 0.8410642299276151, 0.8229028392854474, 0.1104406043887138, 0.3230056762695312, 0.1330508291721344, 0.2336150407791138),
    @("model_9_9_18", 0.9004763783741755, 0.8334162320485992, 0.7856212211623039, 0.8191427976007517, 0.110143318772316, 0.2911181151866913, 0.1794641613960266, 0.2385750710964203),
    @("model_9_9_17", 0.9006999716111493, 0.8342540352995357, 0.7858136035732618, 0.8197878278855568, 0.1098958700895309, 0.2896539568901062, 0.1793031096458435, 0.2377241849899292),
    @("model_9_9_7", 0.9007754824130275, 0.8122664266210528, 0.858666376512736, 0.8261232459451122, 0.1098123043775558, 0.3280790448188782, 0.1183154433965683, 0.2293668836355209),
    @("model_9_9_12", 0.9011186914692234, 0.8372003543634345, 0.7873045064720794, 0.8222993758471486, 0.1094324737787247, 0.2845050394535065, 0.1780550181865692, 0.2344111055135727),
    @("model_9_9_16", 0.9013244183730664, 0.8331146452529903, 0.794917894747006, 0.8217075967247243, 0.1092047840356827, 0.2916451692581177, 0.1716815829277039, 0.235191747546196),
    @("model_9_9_15", 0.9016707838538351, 0.8329530315799208, 0.798728392654051, 0.8227323312273622, 0.1088214591145515, 0.2919275760650635, 0.1684916764497757, 0.2338399887084961),
    @("model_9_9_14", 0.9021213783529796, 0.8365147159007713, 0.7953596983559322, 0.824224136661093, 0.1083227917551994, 0.2857032418251038, 0.1713117361068726, 0.2318720817565918),
    @("model_9_9_13", 0.9022412148782856, 0.8373732141770258, 0.7947864261107062, 0.8246549822070399, 0.1081901788711548, 0.2842029631137848, 0.1717916429042816, 0.2313037663698196),
    @("model_9_9_9", 0.9031982513990655, 0.8281783980683043, 0.8323671560835744, 0.829429406692679, 0.1071310043334961, 0.3002716302871704, 0.1403314620256424, 0.2250056117773056),
    @("model_9_9_10", 0.903285171721601, 0.8314561376012791, 0.8243070313448009, 0.8293211814650086, 0.1070348098874092, 0.2945435047149658, 0.1470788866281509, 0.2251483947038651)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
    }
}
